# LIVEHTA-1904: refresh testdata strings + UI state
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "ICER - ICER RRMM 2022 report - 12/19/2022" population name no longer
# carries the date suffix; both rows that reference the report name need the
# trimmed text.
$ws.Range("D2").Value = "ICER - ICER RRMM 2022 report"
$ws.Range("D4").Value = "ICER - ICER RRMM 2022 report"

# Re-fit the columns whose content widths shifted because of the edit above
# (columns A/B are untouched so they're left alone).
$ws.Columns("C:I").AutoFit() | Out-Null

# Restore the view: scrolled back to the sheet origin with F9 as the active
# selection (previously parked at G7 with C1 pinned to the top-left).
$ws.Range("A1").Select() | Out-Null
$ws.Range("F9").Select() | Out-Null
